# Applies the "DADA2 stuff finished" update to All_F1_tables.xlsx
# - Sheets with Subject "12S"/"16S" (Sheet1,2,4,5,7,8): fill in the
#   previously-zeroed DADA2Spec row (row 5, columns D:H) with the computed
#   Precision/Recall/F1/F0.5/Accuracy values.
# - Sheets with Subject "CO1" (Sheet3,6,9): rename "CO1" -> "COI" for every
#   row (2-15) in column C.

$wb = $excel.ActiveWorkbook

# --- Sheet1 (100 Australian species / 12S) ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D5").Value = 0.9076923076923077
$ws.Range("E5").Value = 0.6344086021505376
$ws.Range("F5").Value = 0.7468354430379748
$ws.Range("G5").Value = 0.8356940509915014
$ws.Range("H5").Value = 0.5959595959595959

# --- Sheet2 (100 Australian species / 16S) ---
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("D5").Value = 0.9803921568627451
$ws.Range("E5").Value = 0.5102040816326531
$ws.Range("F5").Value = 0.6711409395973154
$ws.Range("G5").Value = 0.8278145695364236
$ws.Range("H5").Value = 0.5050505050505051

# --- Sheet3 (100 Australian species / CO1): rename CO1 -> COI ---
$ws = $wb.Worksheets.Item("Sheet3")
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = "COI"
}

# --- Sheet4 (Lutjanidae / 12S) ---
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Range("D5").Value = 0.9333333333333333
$ws.Range("E5").Value = 0.6666666666666666
$ws.Range("F5").Value = 0.7777777777777778
$ws.Range("G5").Value = 0.8641975308641976
$ws.Range("H5").Value = 0.6666666666666666

# --- Sheet5 (Lutjanidae / 16S) ---
$ws = $wb.Worksheets.Item("Sheet5")
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 0.36
$ws.Range("F5").Value = 0.5142857142857143
$ws.Range("G5").Value = 0.6923076923076923
$ws.Range("H5").Value = 0.3703703703703703

# --- Sheet6 (Lutjanidae / CO1): rename CO1 -> COI ---
$ws = $wb.Worksheets.Item("Sheet6")
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = "COI"
}

# --- Sheet7 (Wadjemup / 12S) ---
$ws = $wb.Worksheets.Item("Sheet7")
$ws.Range("D5").Value = 0.95
$ws.Range("E5").Value = 0.5816326530612245
$ws.Range("F5").Value = 0.7215189873417721
$ws.Range("G5").Value = 0.8431952662721893
$ws.Range("H5").Value = 0.5686274509803921

# --- Sheet8 (Wadjemup / 16S) ---
$ws = $wb.Worksheets.Item("Sheet8")
$ws.Range("D5").Value = 0.9285714285714286
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.65
$ws.Range("G5").Value = 0.7926829268292683
$ws.Range("H5").Value = 0.5

# --- Sheet9 (Wadjemup / CO1): rename CO1 -> COI ---
$ws = $wb.Worksheets.Item("Sheet9")
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = "COI"
}
